$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 316, shifting existing rows 316:331 down to 317:332
$ws.Rows("316").Insert()

# Populate the newly inserted row 316 with the new data record
$ws.Cells.Item(316, 1).Value = 6
$ws.Cells.Item(316, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(316, 3).Value = "Metropolitana"
$ws.Cells.Item(316, 4).Value = 44585
$ws.Cells.Item(316, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(316, 5).Value = 13
$ws.Cells.Item(316, 6).Value = 100112032
$ws.Cells.Item(316, 7).Value = "Zapallo italiano"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 400
$ws.Cells.Item(316, 11).Value = 12000
$ws.Cells.Item(316, 12).Value = 13000
$ws.Cells.Item(316, 13).Value = 12575
$ws.Cells.Item(316, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(316, 15).Value = "Región Metropolitana"
$ws.Cells.Item(316, 16).Value = 252
$ws.Cells.Item(316, 17).Value = 50
$ws.Cells.Item(316, 18).Value = "Hortaliza"
